# Actualizacion automatica del mapa (2025-11-28 15:53:00)
#
# Two new case rows were inserted into the "Optical_Power" sheet:
#   - a new row at (final) position 40: caso "-579" / Pedro Rivera 2546
#   - a new row at (final) position 78: caso "7829 " / ALBERDI, JUAN BAUTISTA AV. 2309
# Everything that was at/after each insertion point shifts down by one, which
# is exactly what a real row-Insert does, so we insert a blank row at each
# spot (in top-to-bottom order, so the second insertion point is still
# correct after the first shift) and then fill in the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{
        Row = 40
        A = "-579"; B = "9/2/2025"; C = "Pedro Rivera 2546"; D = "13"; E = "ICD30612785"
        F = "Optical Power"; G = "Pendiente"; H = "Colocar R200 para pedir traspaso de equipos"
        I = 1; J = "Cambio"; K = "Nodo Teco"; L = "Pasante"
        M = -58.462479; N = -34.55765; O = "Saavedra"; P = "Capital Norte"; Q = "COG-K"
        R = "Fuera de Poligono OVL"
    },
    @{
        Row = 78
        A = "7829 "; B = "11/2/2025"; C = "ALBERDI, JUAN BAUTISTA AV. 2309"; D = "7"; E = "810526272"
        F = "Optical Power"; G = "Pendiente"; H = "Picada"
        I = 1; J = "Cambio"; K = "Sin equipos"; L = "Pasante"
        M = -58.460356; N = -34.630793; O = "Boedo"; P = "Capital Sur"; Q = "PCH-J"
        R = "Fuera de Poligono OVL"
    }
)

# Columns A, B, D, E hold values that LOOK numeric/date-like (e.g. "7260",
# "13", "9/2/2025") but are actually stored as plain text in this sheet, so
# those must be written with an explicit text NumberFormat first -
# otherwise Excel would coerce them into numbers/dates and e.g. drop the
# trailing space in "7829 " or turn "9/2/2025" into a date serial.
# The remaining text columns (C, F, G, H, J, K, L, O, P, Q, R) are never
# ambiguous, so they're written as plain text with no format change needed.
# I, M, N are genuine numbers and likewise need no format change.
$forceTextCols = @("A","B","D","E")
$plainTextCols = @("C","F","G","H","J","K","L","O","P","Q","R")
$numCols = @("I","M","N")

foreach ($newRow in $newRows) {
    $r = $newRow.Row

    $ws.Rows.Item($r).Insert()

    foreach ($col in $forceTextCols) {
        $cell = $ws.Range("$col$r")
        $cell.NumberFormat = "@"
        $cell.Value = $newRow[$col]
    }
    foreach ($col in $plainTextCols) {
        $ws.Range("$col$r").Value = $newRow[$col]
    }
    foreach ($col in $numCols) {
        $ws.Range("$col$r").Value = $newRow[$col]
    }
}
